# Apply the "Haste" sheet corrections for the 2025 Fahrplan evaluation:
#  - row 7 ("Nienburg (S)") is removed entirely
#  - row 3 now compares against "Hannover Hbf (RE)" (was "Braunschweig (RE)")
#  - row 6 now compares against "Weetzen (S)" (was "Hannover Hbf (S2)"), which
#    no longer needs the wrapped two-line row height
#  - the active selection on the sheet moves to F3

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Haste")

# Drop the obsolete "Nienburg (S)" comparison row; rows below shift up.
$ws.Rows.Item(7).Delete()

# Row 3: replace "Braunschweig (RE)" figures with "Hannover Hbf (RE)" figures.
$ws.Range("A3").Value = "Hannover Hbf (RE)"
$ws.Range("B3").Value = 18
$ws.Range("C3").Value = 34
$ws.Range("D3").Value = 28
$ws.Range("E3").Value = 29

# Row 6: replace "Hannover Hbf (S2)" figures with "Weetzen (S)" figures.
$ws.Range("A6").Value = "Weetzen (S)"
$ws.Range("B6").Value = 31
$ws.Range("C6").Value = 34
$ws.Range("D6").Value = 26
$ws.Range("E6").Value = 27
$ws.Range("F6").Value = 2

# "Weetzen (S)" is a single short line, so the row no longer needs the
# taller wrapped-text height that "Hannover Hbf (S2)" required.
$ws.Rows.Item(6).AutoFit()

$ws.Range("F3").Select()
